# UF-17724: Fix for PR3 import regarding Ange
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct/replace a batch of duplicated personnummer / phone values on the
# existing rows (2-4): each row should carry its own unique PERSONNR (col K)
# and unique TELEFON (col O) instead of the previously duplicated values.
$ws.Range("K2").Value = "100108-2393"
$ws.Range("O2").Value = "070-1740635"

$ws.Range("O3").Value = "070-1740636"

$ws.Range("O4").Value = "070-1740637"

# K3 (Janne's old placeholder row) gets a distinctly-formatted personnummer
# cell: no border, right aligned + vertically centered, plain number format.
$ws.Range("K3").Value = "250107-2389"
$ws.Range("K3").NumberFormat = "0"
$ws.Range("K3").HorizontalAlignment = -4152
$ws.Range("K3").VerticalAlignment = -4108
$ws.Range("K3").WrapText = $false
$ws.Range("K3").Borders.LineStyle = -4142
$ws.Range("K3").Font.Name = "Arial"
$ws.Range("K3").Font.Size = 10

# --- Append a new row (5) for the new person "JANNE" (Stadsbyggnadsnämnden /
# Ange import) by duplicating the format of row 4 and filling in values.
$ws.Range("A4:AB4").Copy()
$ws.Range("A5:AB5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(5).RowHeight = 45

$ws.Cells.Item(5,1).Value = 101
$ws.Cells.Item(5,2).Value = 622
$ws.Cells.Item(5,3).Value = 0
$ws.Cells.Item(5,4).Value = 1
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = 0
$ws.Cells.Item(5,8).Value = 80704
$ws.Range("I5").Value = "JANSSON"
$ws.Range("J5").Value = "JANNE"
$ws.Range("K5").Value = "000111-2382"
$ws.Range("M5").Value = "VÄGEN 35"
$ws.Range("N5").Value = "123 45 STADEN"
$ws.Range("O5").Value = "070-1740638"
$ws.Range("P5").Value = "08/2619"
$ws.Cells.Item(5,17).Value = 39799
$ws.Cells.Item(5,18).Value = 39799
$ws.Cells.Item(5,19).Value = 40877
$ws.Cells.Item(5,20).Value = 39799
$ws.Cells.Item(5,21).Value = 39799
$ws.Cells.Item(5,22).Value = 39799
$ws.Range("X5").Value = "Stadsbyggnadsnämnden"
$ws.Range("Y5").Value = "Anki Borg"
$ws.Cells.Item(5,28).Value = 1

# K5 gets its own distinct (un-bordered, plain Calibri) number style, locked
# for protection, matching the imported row's origin formatting.
$ws.Range("K5").NumberFormat = "0"
$ws.Range("K5").Borders.LineStyle = -4142
$ws.Range("K5").Font.Name = "Calibri"
$ws.Range("K5").Font.Size = 11
$ws.Range("K5").Locked = $true

# --- Column B (new ID_PERSON-ish column now visible) needs an explicit width.
$ws.Columns.Item(2).ColumnWidth = 11.2

# --- Selection moves to the newly added row.
$ws.Range("K5").Select()

Write-Host "edit applied"
